# culture_collection column (T) is being removed from the MIxS human-skin
# template again (per INSDC2017 review). Deleting the whole column shifts
# all the data/header cells correctly, but comments ("Notes") attached to
# cells are not carried along by the engine's Range.Delete, so we capture
# them first and re-create them afterwards at their new (shifted) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$DeletedColumn = 20  # column T ("culture_collection")

# 1. Snapshot every existing comment: (row, col, text)
$savedComments = @()
foreach ($cm in $ws.Comments) {
    $p = $cm.Parent
    $savedComments += , @($p.Row, $p.Column, $cm.Text())
}

# 2. Remove all comments before shifting cells around, so none are left
#    dangling on the wrong cell once the column shifts. Walk the live
#    collection back-to-front: deleting by a stale forward index/object
#    reference after earlier removals can silently no-op.
for ($i = $ws.Comments.Count; $i -ge 1; $i--) {
    $ws.Comments.Item($i).Delete()
}

# 3. Delete the whole culture_collection column; this shifts the
#    remaining worksheet data/headers left by one column.
$ws.Columns.Item($DeletedColumn).Delete()

# 4. Re-create the comments on their new cells: anything left of the
#    deleted column keeps its place, anything that was on the deleted
#    column itself is dropped, everything to the right moves one column
#    to the left.
foreach ($entry in $savedComments) {
    $row = $entry[0]
    $col = $entry[1]
    $text = $entry[2]

    if ($col -eq $DeletedColumn) {
        continue
    }

    $newCol = $col
    if ($col -gt $DeletedColumn) {
        $newCol = $col - 1
    }

    $target = $ws.Cells.Item($row, $newCol)
    $target.AddComment($text) | Out-Null
}
